# Daily attendance processing - 2025-11-09 13:45:46
# Normalize the "Recorded By" (column G) values on the active sheet by
# sorting the comma-separated list of recorder names/emails for each row
# using an ordinal (code-point) comparison.

function Compare-Ordinal($cmpA, $cmpB) {
    $lenA = $cmpA.Length
    $lenB = $cmpB.Length
    $minLen = $lenA
    if ($lenB -lt $minLen) { $minLen = $lenB }
    $charIdx = 0
    $result = 0
    $found = $false
    while ($charIdx -lt $minLen -and -not $found) {
        $ca = [int][char]$cmpA[$charIdx]
        $cb = [int][char]$cmpB[$charIdx]
        if ($ca -lt $cb) { $result = -1; $found = $true }
        elseif ($ca -gt $cb) { $result = 1; $found = $true }
        $charIdx = $charIdx + 1
    }
    if (-not $found) {
        if ($lenA -lt $lenB) { $result = -1 }
        elseif ($lenA -gt $lenB) { $result = 1 }
        else { $result = 0 }
    }
    return $result
}

function Sort-Ordinal($items) {
    $n = $items.Count
    for ($insIdx = 1; $insIdx -lt $n; $insIdx++) {
        $key = $items[$insIdx]
        $shiftIdx = $insIdx - 1
        while ($shiftIdx -ge 0 -and (Compare-Ordinal $items[$shiftIdx] $key) -gt 0) {
            $items[$shiftIdx + 1] = $items[$shiftIdx]
            $shiftIdx = $shiftIdx - 1
        }
        $items[$shiftIdx + 1] = $key
    }
    return $items
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($rowNum = 2; $rowNum -le $lastRow; $rowNum++) {
    $cell = $ws.Cells.Item($rowNum, 7)
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $rawParts = $value -split ",\s*"
        $parts = @()
        foreach ($p in $rawParts) {
            $parts += $p.Trim()
        }

        $sortedParts = Sort-Ordinal $parts
        $newValue = [string]::Join(", ", $sortedParts)

        if ($newValue -ne $value) {
            $cell.Value = $newValue
        }
    }
}
